$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# D-column values are stored as literal text (prices use "." as a
# thousands separator, e.g. "26.607.08"), so ambiguous numeric-looking
# strings are forced to text via NumberFormat "@" before assignment,
# then ClearFormats() restores the default (General) cell format so
# only the cell VALUE changes, matching the source data.

$ws.Range("D2").Value = '26.607.08'

$ws.Range("D3").Value = '1.718.81'
$ws.Range("E3").Value = '  -1.23%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9979'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.68'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9985'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4923'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.82%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2602'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.93%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06207'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.06%  '

$ws.Range("D10").Value = '1.724.07'
$ws.Range("E10").Value = '  -0.93%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06999'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.75'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6072'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.481'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.51%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.78'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9987'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.12%  '

$ws.Range("D17").Value = '26.445.29'
$ws.Range("E17").Value = '  -0.73%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9980'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007150'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.79%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.36'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.86%  '

$ws.Range("D21").Value = '1.949.00'
$ws.Range("E21").Value = '  -0.94%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.407'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.36%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.511'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.080'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '137.67'
$ws.Range("D25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.29'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.402'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.98%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.743'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.05%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '105.77'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.38%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.917'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.61%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07945'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.640'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04505'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.612'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9990'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.74%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6258'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9359'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.99%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.007'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.46%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.411'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9983'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01507'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.28'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.513'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3839'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.61%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.922'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.80%  '

$ws.Range("E46").Value = '  -2.42%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05370'
$ws.Range("D47").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.739'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.13'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.90%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.44'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.222'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.50%  '
